# Update a handful of transaction/reconciliation figures across three
# sheets. The "Expected Out" total in B1 is a SUM formula and will
# recalculate automatically once its inputs change.

$wb = $excel.ActiveWorkbook

# --- TestRecord sheet -------------------------------------------------
$wsTestRecord = $wb.Worksheets.Item("TestRecord")
$wsTestRecord.Range("A10").Value = 43268
$wsTestRecord.Range("B10").Value = 127.14
$wsTestRecord.Range("E10").Value = "some test textzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzz"

# --- Budget Out sheet ---------------------------------------------------
$wsBudgetOut = $wb.Worksheets.Item("Budget Out")
$wsBudgetOut.Range("C9").Value = 96.22
$wsBudgetOut.Range("F9").Value = "Description007zzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzz"

# --- Expected Out sheet --------------------------------------------------
$wsExpectedOut = $wb.Worksheets.Item("Expected Out")
$wsExpectedOut.Range("B9").Value = 1354.16
$wsExpectedOut.Range("B11").Value = 434.02

# Try to widen the saved window a touch to mirror the author's view state.
# (No-op on hosts that don't expose/persist window geometry.)
$win = $excel.ActiveWindow
$win.Height = 10800
